$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("D2").Value = "2016-01-13 04:45:33"
$ws.Range("G2").Value = "2016-01-13 04:46:45"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("D2").Value = "2016-01-13 04:45:54"
$ws2.Range("G2").Value = "2016-01-13 04:47:19"
